{"js": "// Renumber the \"\u65e5\u5fd7\u914d\u7f6e\" heading from \"3.1.3\" to \"3.1.4\" (a new \"3.1.3 \u4ee3\u7801\u751f\u6210\u914d\u7f6e\"\n// section was inserted earlier in the document, so the logging-configuration\n// heading that used to be the 3rd sub-section is now the 4th).\n//\n// Before: \"3.1.3 \u65e5\u5fd7\u914d\u7f6e\"\n// After:  \"3.1.4 \u65e5\u5fd7\u914d\u7f6e\"\n//\n// The document also carries a \"_GoBack\" bookmark (Word's \"last edit location\"\n// marker) that, after the real edit, ends up sitting right after the new \"4\"\n// digit (i.e. between \"3.1.4\" and the following space + \"\u65e5\u5fd7\u914d\u7f6e\"), so we\n// reproduce that placement too.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the unique heading paragraph whose full text is \"3.1.3 \u65e5\u5fd7\u914d\u7f6e\".\n// (There is an earlier, unrelated \"3.1.3 \u4ee3\u7801\u751f\u6210\u914d\u7f6e\" heading, so we match on\n// the full heading text rather than just the \"3.1.3 \" prefix.)\nconst headingText = \"3.1.3 \\u65E5\\u5FD7\\u914D\\u7F6E\"; // \"3.1.3 \u65e5\u5fd7\u914d\u7f6e\"\nlet heading = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === headingText) {\n    heading = paragraphs.items[i];\n    break;\n  }\n}\n\nif (heading) {\n  // Within that paragraph, locate the section-number digit that needs to\n  // change: the second \"3\" in \"3.1.3\" (the first \"3\" is the chapter number\n  // and must stay untouched).\n  const digitMatches = heading.search(\"3\", { matchCase: true });\n  digitMatches.load(\"items\");\n  await context.sync();\n\n  const sectionDigit = digitMatches.items[digitMatches.items.length - 1];\n  sectionDigit.insertText(\"4\", Word.InsertLocation.replace);\n  await context.sync();\n\n  // Move the \"_GoBack\" bookmark so it again sits right after the edited\n  // digit (matching where Word leaves it after an in-place text edit).\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n\n  const fourMatches = heading.search(\"4\", { matchCase: true });\n  fourMatches.load(\"items\");\n  await context.sync();\n\n  const fourRange = fourMatches.items[fourMatches.items.length - 1];\n  const afterFour = fourRange.getRange(Word.RangeLocation.after);\n  afterFour.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Renumber the \"\u65e5\u5fd7\u914d\u7f6e\" heading from \"3.1.3\" to \"3.1.4\" (a new \"3.1.3 \u4ee3\u7801\u751f\u6210\u914d\u7f6e\"\n# section was inserted earlier in the document, so the logging-configuration\n# heading that used to be the 3rd sub-section is now the 4th).\n#\n# Before: \"3.1.3 \u65e5\u5fd7\u914d\u7f6e\"\n# After:  \"3.1.4 \u65e5\u5fd7\u914d\u7f6e\"\n#\n# The document also carries a \"_GoBack\" bookmark (Word's \"last edit location\"\n# marker) that, after the real edit, ends up sitting right after the new \"4\"\n# digit (i.e. between \"3.1.4\" and the following space + \"\u65e5\u5fd7\u914d\u7f6e\"), so we\n# reproduce that placement too.\n\n$d = $word.ActiveDocument\n\n# Locate the unique heading \"3.1.3 \u65e5\u5fd7\u914d\u7f6e\" (there is an earlier, unrelated\n# \"3.1.3 \u4ee3\u7801\u751f\u6210\u914d\u7f6e\" heading, so search on the full heading text).\n$needle = \"3.1.3 \" + [char]0x65E5 + [char]0x5FD7 + [char]0x914D + [char]0x7F6E\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute($needle)\n\nif ($found) {\n    $headingStart = $rng.Start\n\n    # Within \"3.1.3 \", the section-number digit to bump is the second \"3\"\n    # (4 characters in: \"3\",\".\",\"1\",\".\", -> index 4 is the second \"3\").\n    $digitRange = $d.Range($headingStart + 4, $headingStart + 5)\n    if ($digitRange.Text -eq \"3\") {\n        $digitRange.Text = \"4\"\n    }\n\n    # Re-anchor the \"_GoBack\" bookmark immediately after the new \"4\" digit,\n    # matching where Word leaves it after an in-place text edit.\n    if ($d.Bookmarks.Exists(\"_GoBack\")) {\n        $d.Bookmarks(\"_GoBack\").Delete()\n    }\n    $bookmarkPos = $d.Range($headingStart + 5, $headingStart + 5)\n    $d.Bookmarks.Add(\"_GoBack\", $bookmarkPos)\n}\n"}
